$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K4").Value = 13
$ws.Range("L4").Value = 21.23

$ws.Range("K5").Value = 31.26
$ws.Range("L5").Value = 45.66

$ws.Range("K6").Value = 35
$ws.Range("L6").Value = 92.09

$ws.Range("K7").Value = 161.113
$ws.Range("L7").Value = 121.38

$ws.Range("K8").Value = 245.373
$ws.Range("L8").Value = 291.37
